$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.Formula = "=" + [char]34 + $value + [char]34
    $r.Copy()
    $r.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
}

$excel.CutCopyMode = $false

Set-TextValue "D2" "26.786.15"
Set-TextValue "E2" "  +0.56%  "
Set-TextValue "D3" "1.644.82"
Set-TextValue "E3" "  +0.15%  "
Set-TextValue "E4" "  +0.25%  "
Set-TextValue "D5" "216.55"
Set-TextValue "E5" "  +0.56%  "
Set-TextValue "D6" "0.501"
Set-TextValue "E6" "  -0.70%  "
Set-TextValue "E7" "  +0.29%  "
Set-TextValue "E9" "  +0.21%  "
Set-TextValue "E10" "  -0.30%  "
Set-TextValue "D11" "0.0843"
Set-TextValue "E11" "  +0.15%  "
Set-TextValue "D12" "1.659.37"
Set-TextValue "E12" "  +1.43%  "
Set-TextValue "D13" "4.16"
Set-TextValue "E13" "  -0.85%  "
Set-TextValue "E14" "  -0.41%  "
Set-TextValue "E15" "  -1.17%  "
Set-TextValue "D16" "26.772.39"
Set-TextValue "E16" "  +0.40%  "
Set-TextValue "D17" "0.0₃0736"
Set-TextValue "D18" "214.27"
Set-TextValue "E18" "  -1.11%  "
Set-TextValue "E19" "  +0.28%  "
Set-TextValue "E20" "  +0.96%  "
Set-TextValue "E21" "  +12.36%  "
Set-TextValue "E22" "  -0.82%  "
Set-TextValue "E23" "  -1.57%  "
Set-TextValue "D24" "146.76"
Set-TextValue "E24" "  +0.82%  "
Set-TextValue "E25" "  +0.25%  "
Set-TextValue "E26" "  -0.83%  "
Set-TextValue "D27" "7.17"
Set-TextValue "E27" "  +0.20%  "
Set-TextValue "D28" "15.65"
Set-TextValue "E28" "  -0.75%  "
Set-TextValue "E29" "  -1.52%  "
Set-TextValue "E30" "  +0.68%  "
Set-TextValue "E31" "  -0.46%  "
Set-TextValue "E32" "  -1.26%  "
Set-TextValue "D33" "1.291.58"
Set-TextValue "E33" "  +0.93%  "
Set-TextValue "E34" "  -0.45%  "
Set-TextValue "E35" "  +1.49%  "
Set-TextValue "E36" "  -2.52%  "
Set-TextValue "E37" "  +0.94%  "
Set-TextValue "E38" "  -0.64%  "
Set-TextValue "E39" "  +0.23%  "
Set-TextValue "E40" "  -1.25%  "
Set-TextValue "E41" "  -0.83%  "
Set-TextValue "E42" "  -2.72%  "
Set-TextValue "D43" "1.783.41"
Set-TextValue "E43" "  +0.11%  "
Set-TextValue "D44" "61.83"
Set-TextValue "E44" "  +3.11%  "
Set-TextValue "D45" "91.58"
Set-TextValue "E45" "  -0.50%  "
Set-TextValue "E46" "  +0.89%  "
Set-TextValue "E47" "  -1.08%  "
Set-TextValue "E48" "  +0.85%  "
Set-TextValue "D49" "7.65"
Set-TextValue "E49" "  -1.62%  "
Set-TextValue "D50" "0.0971"
Set-TextValue "E50" "  +0.12%  "
Set-TextValue "D51" "0.406"
Set-TextValue "E51" "  -0.04%  "

$excel.CutCopyMode = $false
